$d = $word.ActiveDocument

# 1) " on May 01, 2022." -> " on May 05, 2022." (unique in doc)
$d.Content.Find.Execute(" on May 01, 2022.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " on May 05, 2022.", 2)

# 2) " license is suspended from May 01, 2022" -> " license is suspended from May 05, 2022" (unique in doc)
$d.Content.Find.Execute(" license is suspended from May 01, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, " license is suspended from May 05, 2022", 2)

# 3) "June 30, 2022" -> "July 04, 2022" (unique in doc)
$d.Content.Find.Execute("June 30, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "July 04, 2022", 2)

# 4) remaining standalone "May 01, 2022" (in "...in full by May 01, 2022.") -> "May 05, 2022"
# After replacements 1 and 2, this is now the only remaining match in the document.
$d.Content.Find.Execute("May 01, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "May 05, 2022", 2)
